$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '69.515.74'
Set-TextValue 2 5 '  +6.00%  '
Set-TextValue 3 4 '3.572.84'
Set-TextValue 3 5 '  +5.18%  '
Set-TextValue 4 5 '  -0.06%  '
Set-TextValue 5 4 '591.86'
Set-TextValue 5 5 '  +5.58%  '
Set-TextValue 6 4 '192.21'
Set-TextValue 6 5 '  +9.27%  '
Set-TextValue 7 4 '0.642'
Set-TextValue 7 5 '  +1.59%  '
Set-TextValue 8 4 '3.569.12'
Set-TextValue 8 5 '  +5.28%  '
Set-TextValue 9 4 '0.999'
Set-TextValue 9 5 '  -0.08%  '
Set-TextValue 10 5 '  +5.57%  '
Set-TextValue 11 4 '0.661'
Set-TextValue 11 5 '  +3.97%  '
Set-TextValue 12 4 '58.14'
Set-TextValue 12 5 '  +8.70%  '
Set-TextValue 13 4 '0.0000291'
Set-TextValue 13 5 '  +4.96%  '
Set-TextValue 14 5 '  +5.26%  '
Set-TextValue 15 4 '4.139.33'
Set-TextValue 15 5 '  +5.07%  '
Set-TextValue 16 4 '19.30'
Set-TextValue 16 5 '  +5.61%  '
Set-TextValue 17 4 '3.568.16'
Set-TextValue 17 5 '  +5.63%  '
Set-TextValue 18 4 '69.420.47'
Set-TextValue 18 5 '  +6.01%  '
Set-TextValue 19 4 '12.44'
Set-TextValue 19 5 '  +5.03%  '
Set-TextValue 20 5 '  +0.60%  '
Set-TextValue 21 5 '  +4.65%  '
Set-TextValue 22 4 '502.80'
Set-TextValue 22 5 '  +5.05%  '
Set-TextValue 23 4 '5.51'
Set-TextValue 23 5 '  +11.80%  '
Set-TextValue 24 4 '17.15'
Set-TextValue 24 5 '  +19.92%  '
Set-TextValue 25 5 '  +8.12%  '
Set-TextValue 26 4 '91.13'
Set-TextValue 26 5 '  +1.81%  '
Set-TextValue 27 5 '  +4.77%  '
Set-TextValue 28 4 '11.17'
Set-TextValue 28 5 '  +5.01%  '
Set-TextValue 29 4 '9.33'
Set-TextValue 29 5 '  +6.98%  '
Set-TextValue 30 4 '32.08'
Set-TextValue 30 5 '  +2.69%  '
Set-TextValue 31 4 '7.56'
Set-TextValue 31 5 '  +15.27%  '
$ws.Cells.Item(32, 2).Value = 'Cosmos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 32 4 '12.12'
Set-TextValue 32 5 '  +5.40%  '
$ws.Cells.Item(33, 2).Value = 'Bittensor'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 33 4 '615.35'
Set-TextValue 33 5 '  +6.86%  '
Set-TextValue 34 4 '65.49'
Set-TextValue 34 5 '  +4.06%  '
Set-TextValue 35 5 '  +6.60%  '
Set-TextValue 36 4 '0.0₃0836'
Set-TextValue 36 5 '  +12.88%  '
Set-TextValue 37 5 '  +4.61%  '
$ws.Cells.Item(38, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 38 4 '38.09'
Set-TextValue 38 5 '  +6.29%  '
$ws.Cells.Item(39, 2).Value = 'Dai'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 39 4 '1.00'
Set-TextValue 39 5 '  -0.05%  '
Set-TextValue 40 5 '  +6.27%  '
Set-TextValue 41 4 '3.61'
Set-TextValue 41 5 '  -1.44%  '
Set-TextValue 42 4 '3.322.03'
Set-TextValue 42 5 '  +7.47%  '
Set-TextValue 43 4 '3.12'
Set-TextValue 43 5 '  +11.49%  '
Set-TextValue 44 4 '2.71'
Set-TextValue 44 5 '  +11.44%  '
Set-TextValue 45 4 '0.0441'
Set-TextValue 45 5 '  +5.66%  '
Set-TextValue 46 4 '2.89'
Set-TextValue 46 5 '  +19.91%  '
Set-TextValue 47 4 '3.32'
Set-TextValue 47 5 '  +4.65%  '
Set-TextValue 48 5 '  +2.33%  '
Set-TextValue 49 4 '9.12'
Set-TextValue 49 5 '  +8.22%  '
Set-TextValue 50 5 '  +4.45%  '
Set-TextValue 51 5 '  -0.03%  '
